$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Range("D2").Value = 119097
$ws.Range("E2").Value = 3953
$ws.Range("F2").Value = 3953
$ws.Range("G2").Value = -4417
$ws.Range("H2").Value = -6129
$ws.Range("I2").Value = -6354
$ws.Range("J2").Value = 224
$ws.Range("K2").Value = 234657
$ws.Range("L2").Value = 212646
$ws.Range("M2").Value = 22012
$ws.Range("N2").Value = 20909
$ws.Range("O2").Value = 1102
$ws.Range("P2").Value = 2989
$ws.Range("Q2").Value = 18636
$ws.Range("R2").Value = -10649
$ws.Range("S2").Value = -11379
$ws.Range("T2").Value = 11206
$ws.Range("U2").Value = 7429
$ws.Range("V2").Value = 149427
$ws.Range("W2").Value = 3.32
$ws.Range("X2").Value = -5.15
$ws.Range("Y2").Value = -26.82
$ws.Range("Z2").Value = -2.64
$ws.Range("AA2").Value = 966.0599999999999
$ws.Range("AB2").Value = 523.3200000000001
$ws.Range("AC2").Value = -9471
$ws.Range("AD2").Value = -4.47
$ws.Range("AE2").Value = 31170
$ws.Range("AF2").Value = 1.36
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 65970726
$ws.Range("D3").Value = 115448
$ws.Range("E3").Value = 8831
$ws.Range("F3").Value = 8831
$ws.Range("G3").Value = -4864
$ws.Range("H3").Value = -5630
$ws.Range("I3").Value = -5650
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 241804
$ws.Range("L3").Value = 216813
$ws.Range("M3").Value = 24990
$ws.Range("N3").Value = 23871
$ws.Range("O3").Value = 1119
$ws.Range("P3").Value = 3698
$ws.Range("Q3").Value = 27280
$ws.Range("R3").Value = 4187
$ws.Range("S3").Value = -29966
$ws.Range("T3").Value = 17427
$ws.Range("U3").Value = 9853
$ws.Range("V3").Value = 149870
$ws.Range("W3").Value = 7.65
$ws.Range("X3").Value = -4.88
$ws.Range("Y3").Value = -25.23
$ws.Range("Z3").Value = -2.36
$ws.Range("AA3").Value = 867.59
$ws.Range("AB3").Value = 377.85
$ws.Range("AC3").Value = -7454
$ws.Range("AD3").Value = -3.49
$ws.Range("AE3").Value = 30304
$ws.Range("AF3").Value = 0.86
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 77661527
$ws.Range("D4").Value = 117319
$ws.Range("E4").Value = 11208
$ws.Range("F4").Value = 11208
$ws.Range("G4").Value = -7174
$ws.Range("H4").Value = -5568
$ws.Range("I4").Value = -5649
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 239565
$ws.Range("L4").Value = 220822
$ws.Range("M4").Value = 18744
$ws.Range("N4").Value = 17607
$ws.Range("O4").Value = 1137
$ws.Range("P4").Value = 3698
$ws.Range("Q4").Value = 28063
$ws.Range("R4").Value = -8735
$ws.Range("S4").Value = -18293
$ws.Range("T4").Value = 11451
$ws.Range("U4").Value = 16611
$ws.Range("V4").Value = 145529
$ws.Range("W4").Value = 9.550000000000001
$ws.Range("X4").Value = -4.75
$ws.Range("Y4").Value = -27.24
$ws.Range("Z4").Value = -2.31
$ws.Range("AA4").Value = 1178.13
$ws.Range("AB4").Value = 110.88
$ws.Range("AC4").Value = -7171
$ws.Range("AD4").Value = -3.59
$ws.Range("AE4").Value = 22352
$ws.Range("AF4").Value = 1.15
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 77661527
$ws.Range("D5").Value = 120922
$ws.Range("E5").Value = 9398
$ws.Range("F5").Value = 9398
$ws.Range("G5").Value = 11217
$ws.Range("H5").Value = 8019
$ws.Range("I5").Value = 7915
$ws.Range("J5").Value = 104
$ws.Range("K5").Value = 246487
$ws.Range("L5").Value = 208976
$ws.Range("M5").Value = 37511
$ws.Range("N5").Value = 36294
$ws.Range("O5").Value = 1217
$ws.Range("P5").Value = 4798
$ws.Range("Q5").Value = 28068
$ws.Range("R5").Value = -20412
$ws.Range("S5").Value = -10362
$ws.Range("T5").Value = 18765
$ws.Range("U5").Value = 9302
$ws.Range("V5").Value = 135847
$ws.Range("W5").Value = 7.77
$ws.Range("X5").Value = 6.63
$ws.Range("Y5").Value = 29.37
$ws.Range("Z5").Value = 3.3
$ws.Range("AA5").Value = 557.11
$ws.Range("AB5").Value = 317.41
$ws.Range("AC5").Value = 8631
$ws.Range("AD5").Value = 3.92
$ws.Range("AE5").Value = 37824
$ws.Range("AF5").Value = 0.89
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 0.74
$ws.Range("AI5").Value = 3.04
$ws.Range("AJ5").Value = 94844634
$ws.Range("D6").Value = 130203
$ws.Range("E6").Value = 6403
$ws.Range("F6").Value = 6403
$ws.Range("G6").Value = -2087
$ws.Range("H6").Value = -1857
$ws.Range("I6").Value = -1930
$ws.Range("K6").Value = 255797
$ws.Range("L6").Value = 225479
$ws.Range("M6").Value = 30318
$ws.Range("N6").Value = 29087
$ws.Range("P6").Value = 4798
$ws.Range("Q6").Value = 27926
$ws.Range("R6").Value = -6579
$ws.Range("S6").Value = -14075
$ws.Range("T6").Value = 12762
$ws.Range("U6").Value = 15164
$ws.Range("V6").Value = 145348
$ws.Range("W6").Value = 4.92
$ws.Range("X6").Value = -1.43
$ws.Range("Y6").Value = -5.91
$ws.Range("Z6").Value = -0.74
$ws.Range("AA6").Value = 743.72
$ws.Range("AB6").Value = 214.26
$ws.Range("AC6").Value = -2012
$ws.Range("AD6").Value = -16.43
$ws.Range("AE6").Value = 30313
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 0.76
$ws.Range("AI6").Value = -12.46
$ws.Range("AJ6").Value = 94844634
$ws.Range("D7").Value = 128246
$ws.Range("E7").Value = 1432
$ws.Range("G7").Value = -8259
$ws.Range("H7").Value = -6210
$ws.Range("I7").Value = -6207
$ws.Range("K7").Value = 263716
$ws.Range("L7").Value = 237450
$ws.Range("M7").Value = 26267
$ws.Range("N7").Value = 24800
$ws.Range("P7").Value = 4799
$ws.Range("Q7").Value = 20572
$ws.Range("R7").Value = -18458
$ws.Range("S7").Value = -3645
$ws.Range("T7").Value = 16735
$ws.Range("U7").Value = 6145
$ws.Range("W7").Value = 1.12
$ws.Range("X7").Value = -4.84
$ws.Range("Y7").Value = -23.04
$ws.Range("Z7").Value = -2.39
$ws.Range("AA7").Value = 904
$ws.Range("AC7").Value = -6469
$ws.Range("AD7").Value = -3.66
$ws.Range("AE7").Value = 25845
$ws.Range("AF7").Value = 0.92
$ws.Range("AG7").Value = 203
$ws.Range("AH7").Value = 0.86
$ws.Range("AI7").Value = -3.1
$ws.Range("D8").Value = 132037
$ws.Range("E8").Value = 5018
$ws.Range("G8").Value = -52
$ws.Range("H8").Value = 58
$ws.Range("I8").Value = 43
$ws.Range("K8").Value = 260291
$ws.Range("L8").Value = 233767
$ws.Range("M8").Value = 26524
$ws.Range("N8").Value = 24865
$ws.Range("P8").Value = 4864
$ws.Range("Q8").Value = 24641
$ws.Range("R8").Value = -13339
$ws.Range("S8").Value = -8182
$ws.Range("T8").Value = 12661
$ws.Range("U8").Value = 12247
$ws.Range("W8").Value = 3.8
$ws.Range("X8").Value = 0.04
$ws.Range("Y8").Value = 0.17
$ws.Range("Z8").Value = 0.02
$ws.Range("AA8").Value = 881.34
$ws.Range("AC8").Value = 45
$ws.Range("AD8").Value = 528.48
$ws.Range("AE8").Value = 25913
$ws.Range("AF8").Value = 0.91
$ws.Range("AG8").Value = 209
$ws.Range("AH8").Value = 0.88
$ws.Range("AI8").Value = 461.23
$ws.Range("D9").Value = 136554
$ws.Range("E9").Value = 7190
$ws.Range("G9").Value = 2110
$ws.Range("H9").Value = 1414
$ws.Range("I9").Value = 1405
$ws.Range("K9").Value = 259571
$ws.Range("L9").Value = 231538
$ws.Range("M9").Value = 28032
$ws.Range("N9").Value = 26195
$ws.Range("P9").Value = 4929
$ws.Range("Q9").Value = 25249
$ws.Range("R9").Value = -15672
$ws.Range("S9").Value = -6505
$ws.Range("T9").Value = 14843
$ws.Range("U9").Value = 9556
$ws.Range("W9").Value = 5.27
$ws.Range("X9").Value = 1.04
$ws.Range("Y9").Value = 5.5
$ws.Range("Z9").Value = 0.54
$ws.Range("AA9").Value = 825.97
$ws.Range("AC9").Value = 1464
$ws.Range("AD9").Value = 16.15
$ws.Range("AE9").Value = 27299
$ws.Range("AF9").Value = 0.87
$ws.Range("AG9").Value = 209
$ws.Range("AH9").Value = 0.88
$ws.Range("AI9").Value = 14.1
